# Update "F" column (想去人数 / interested-count) values across all sheets
# to match the data snapshot regenerated at commit 7921097.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1688
$ws.Range("F6").Value = 564
$ws.Range("F7").Value = 261
$ws.Range("F8").Value = 1227
$ws.Range("F9").Value = 984
$ws.Range("F12").Value = 699
$ws.Range("F15").Value = 772
$ws.Range("F18").Value = 1319
$ws.Range("F19").Value = 28
$ws.Range("F20").Value = 5
$ws.Range("F21").Value = 52
$ws.Range("F22").Value = 99
$ws.Range("F23").Value = 20
$ws.Range("F24").Value = 1234
$ws.Range("F25").Value = 313
$ws.Range("F26").Value = 421
$ws.Range("F27").Value = 118
$ws.Range("F28").Value = 83
$ws.Range("F29").Value = 2505
$ws.Range("F37").Value = 48
$ws.Range("F38").Value = 887
$ws.Range("F39").Value = 24

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 47
$ws.Range("F10").Value = 371
$ws.Range("F12").Value = 8
$ws.Range("F13").Value = 203
$ws.Range("F15").Value = 247
$ws.Range("F18").Value = 47
$ws.Range("F20").Value = 604
$ws.Range("F22").Value = 14
$ws.Range("F23").Value = 427
$ws.Range("F26").Value = 183
$ws.Range("F28").Value = 177

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2021
$ws.Range("F6").Value = 2146
$ws.Range("F7").Value = 843
$ws.Range("F8").Value = 803
$ws.Range("F11").Value = 838
$ws.Range("F12").Value = 130

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2021
$ws.Range("F5").Value = 2146
$ws.Range("F6").Value = 1688
$ws.Range("F8").Value = 843
$ws.Range("F9").Value = 803
$ws.Range("F13").Value = 564
$ws.Range("F14").Value = 261
$ws.Range("F16").Value = 1227
$ws.Range("F17").Value = 984
$ws.Range("F18").Value = 838
$ws.Range("F21").Value = 130
$ws.Range("F24").Value = 772
$ws.Range("F26").Value = 28
$ws.Range("F27").Value = 52
$ws.Range("F28").Value = 99
$ws.Range("F29").Value = 20
$ws.Range("F30").Value = 1234
$ws.Range("F31").Value = 313
$ws.Range("F32").Value = 421
$ws.Range("F33").Value = 83
$ws.Range("F34").Value = 2505
$ws.Range("F41").Value = 887
$ws.Range("F42").Value = 47
$ws.Range("F43").Value = 14
$ws.Range("F45").Value = 183
$ws.Range("F46").Value = 177

